# Update forecast summary workbook with corrected forecast output.
#
# 1. "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
#    (B) between "Week" (A) and "ASIN" (old B, now C); fill it with the
#    week-start dates as plain text. Shorten the "Week" labels from
#    "W01".."W09" to "W1".."W9" (W10-W16 stay the same). Fix two
#    MyForecast values (now column D) that were wrong. Turn the
#    "is_holiday_week" column (now J) into real boolean cells.
# 2. "Summary" sheet: bump "Total Forecast (16 Weeks)" from 97 to 98.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert the new Week_Start_Date column at B, shifting ASIN.. right ---
$ws.Columns.Item(2).EntireColumn.Insert()

$ws.Range("B1").Value = "Week_Start_Date"

# Week-start date for each data row (r=2..17). Force text formatting
# first so Excel doesn't auto-convert the ISO-looking strings into date
# serial numbers.
$weekStartDates = @{
    2  = "2025-01-05"
    3  = "2025-01-12"
    4  = "2025-01-19"
    5  = "2025-01-26"
    6  = "2025-02-02"
    7  = "2025-02-09"
    8  = "2025-02-16"
    9  = "2025-02-23"
    10 = "2025-03-02"
    11 = "2025-03-09"
    12 = "2025-03-16"
    13 = "2025-03-23"
    14 = "2025-03-30"
    15 = "2025-04-06"
    16 = "2025-04-13"
    17 = "2025-04-20"
}

$ws.Range("B2:B17").NumberFormat = "@"
foreach ($r in $weekStartDates.Keys) {
    $ws.Cells.Item($r, 2).Value = $weekStartDates[$r]
}

# --- Shorten the Week labels in column A: "W01" -> "W1", etc. ---
$weekLabels = @{
    2  = "W1"
    3  = "W2"
    4  = "W3"
    5  = "W4"
    6  = "W5"
    7  = "W6"
    8  = "W7"
    9  = "W8"
    10 = "W9"
}
foreach ($r in $weekLabels.Keys) {
    $ws.Cells.Item($r, 1).Value = $weekLabels[$r]
}

# --- Correct two MyForecast values (column D after the column insert) ---
$ws.Range("D9").Value = 7
$ws.Range("D16").Value = 6

# --- Make the is_holiday_week column (now J) hold real booleans ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 10).Value = $false
}

# --- Summary sheet: Total Forecast (16 Weeks) 97 -> 98 ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "98"
